$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet: BothFilter - insert the new "Remove years after 1999..." filtering
# step into both the northern and southern filtering blocks, and relabel the
# Stephens-MacCall filter rows to indicate which filter is selected vs the
# alternate one.
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("BothFilter")

# Insert a new blank row above the current row 10 (the "NA" separator that
# starts the southern filtering block), shifting the old rows 10-14 down to
# rows 11-15.
$ws.Rows.Item(10).Insert()

# New row 10: "Remove years after 1999..." filtering step for the northern
# block. Written in this order so the new shared strings line up as
# 38 ("Remove years..."), 39 ("544"), 40 ("220").
$ws.Range("A10").Value = "Remove years after 1999 due to regulation changes and with fewer than 20 trips"
$ws.Range("B10").NumberFormat = "@"
$ws.Range("B10").Value = "544"
$ws.Range("C10").NumberFormat = "@"
$ws.Range("C10").Value = "220"

# Relabel the (still at rows 8 & 9, unaffected by the insert above row 10)
# Stephens-MacCall filter rows for the northern block. New shared strings
# 41 ("...selected filter)") then 42 ("Alternate...").
$ws.Range("A8").Value = "Stephens-MacCall filter (keep all positives - selected filter)"
$ws.Range("A9").Value = "Alternate Stephens-MacCall filter (keep only above threshold)"

# Same relabelling for the shifted southern-block rows (now 14 & 15).
$ws.Range("A14").Value = "Stephens-MacCall filter (keep all positives - selected filter)"
$ws.Range("A15").Value = "Alternate Stephens-MacCall filter (keep only above threshold)"

# Append a new row 16: the southern-block counterpart of the new row 10.
$ws.Range("A16").Value = "Remove years after 1999 due to regulation changes and with fewer than 20 trips"
$ws.Range("B16").Value = 475
$ws.Range("C16").Value = 342
$ws.Range("C16").NumberFormat = "@"

# Widen column A to (roughly) fit the new, longer filter-description text.
$ws.Columns.Item(1).ColumnWidth = 65

# Leave the selection on the newly added last cell, like the source edit.
[void]$ws.Range("C16").Select()
